$wb = $excel.ActiveWorkbook

# Sheet "OFF" - row 2 updates
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 190
$wsOff.Range("C2").Value = 143
$wsOff.Range("D2").Value = 36
$wsOff.Range("E2").Value = 16
$wsOff.Range("F2").Value = 2

# Sheet "DEF" - row 2 updates
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 139
$wsDef.Range("C2").Value = 104
$wsDef.Range("D2").Value = 32
$wsDef.Range("E2").Value = 17
